$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updates: row -> (C new value, E new value)
$updates = @(
    @{Row=13;  C=187865;   E=1168297054},
    @{Row=48;  C=150637;   E=275745293},
    @{Row=81;  C=88359;    E=499734008},
    @{Row=91;  C=18886;    E=75390235},
    @{Row=121; C=1306429;  E=2275588183},
    @{Row=129; C=633892;   E=3436494215},
    @{Row=132; C=586065;   E=3473732208},
    @{Row=144; C=25088;    E=92615653},
    @{Row=186; C=236842;   E=1190190297},
    @{Row=207; C=154665;   E=753774902},
    @{Row=240; C=205941;   E=1070112032},
    @{Row=246; C=18838;    E=71629583}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
